# Update cryptos list (simulated daily data refresh from GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "44.757.26"
$ws.Cells.Item(2, 5).Value = "  -1.94%  "

$ws.Cells.Item(3, 4).Value = "2.339.61"
$ws.Cells.Item(3, 5).Value = "  -2.45%  "

$ws.Cells.Item(4, 5).Value = "  +0.23%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "327.34"
$ws.Cells.Item(5, 5).Value = "  +2.29%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "105.20"
$ws.Cells.Item(6, 5).Value = "  -8.97%  "

$ws.Cells.Item(7, 5).Value = "  -1.29%  "

$ws.Cells.Item(8, 5).Value = "  +0.10%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.608"
$ws.Cells.Item(9, 5).Value = "  -3.42%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "40.25"
$ws.Cells.Item(10, 5).Value = "  -5.82%  "

$ws.Cells.Item(11, 5).Value = "  -2.24%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "8.32"
$ws.Cells.Item(12, 5).Value = "  -4.37%  "

$ws.Cells.Item(13, 5).Value = "  -1.32%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.963"

$ws.Cells.Item(15, 4).Value = "2.703.58"

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "15.17"
$ws.Cells.Item(16, 5).Value = "  -5.29%  "

$ws.Cells.Item(17, 4).Value = "2.335.96"
$ws.Cells.Item(17, 5).Value = "  -2.83%  "

$ws.Cells.Item(18, 4).Value = "44.875.15"
$ws.Cells.Item(18, 5).Value = "  -1.70%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "15.17"
$ws.Cells.Item(19, 5).Value = "  +10.73%  "

$ws.Cells.Item(20, 5).Value = "  -4.06%  "

$ws.Cells.Item(21, 5).Value = "  -2.97%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "3.64"
$ws.Cells.Item(22, 5).Value = "  +1.14%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "72.58"
$ws.Cells.Item(23, 5).Value = "  -3.42%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "255.70"
$ws.Cells.Item(24, 5).Value = "  -3.70%  "

$ws.Cells.Item(25, 5).Value = "  -4.63%  "

$ws.Cells.Item(26, 5).Value = "  +0.11%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "11.16"
$ws.Cells.Item(27, 5).Value = "  -2.21%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "7.31"
$ws.Cells.Item(28, 5).Value = "  -6.36%  "

$ws.Cells.Item(29, 5).Value = "  -3.18%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.0948"
$ws.Cells.Item(30, 5).Value = "  -4.58%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "21.87"
$ws.Cells.Item(31, 5).Value = "  -4.29%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "36.44"
$ws.Cells.Item(32, 5).Value = "  -9.22%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "165.98"
$ws.Cells.Item(33, 5).Value = "  -4.17%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.76"
$ws.Cells.Item(34, 5).Value = "  -6.11%  "

$ws.Cells.Item(35, 5).Value = "  -2.17%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "3.19"
$ws.Cells.Item(36, 5).Value = "  +2.77%  "

$ws.Cells.Item(37, 5).Value = "  -2.74%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "4.66"
$ws.Cells.Item(38, 5).Value = "  -7.45%  "

$ws.Cells.Item(39, 5).Value = "  +7.00%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.91"
$ws.Cells.Item(40, 5).Value = "  -7.08%  "

$ws.Cells.Item(41, 5).Value = "  -4.60%  "

$ws.Cells.Item(42, 4).Value = "1.869.93"
$ws.Cells.Item(42, 5).Value = "  +13.32%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "95.03"
$ws.Cells.Item(43, 5).Value = "  -5.25%  "

$ws.Cells.Item(44, 2).Value = "THORChain"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "6.01"
$ws.Cells.Item(44, 5).Value = "  +2.45%  "

$ws.Cells.Item(45, 5).Value = "  +0.10%  "

$ws.Cells.Item(46, 2).Value = "MultiversX"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "68.29"
$ws.Cells.Item(46, 5).Value = "  -5.59%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.224"
$ws.Cells.Item(47, 5).Value = "  -7.70%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "12.71"
$ws.Cells.Item(48, 5).Value = "  -7.64%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "83.03"
$ws.Cells.Item(49, 5).Value = "  -7.45%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "9.22"
$ws.Cells.Item(50, 5).Value = "  -2.76%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "110.34"
$ws.Cells.Item(51, 5).Value = "  -5.20%  "

